$d = $word.ActiveDocument

# Extend the final sentence of the document: "...Windalya va être attaquée."
# becomes "...Windalya va être attaquée, et décide de rentrer rapidos à
# Windalya, Kris avec lui. Katar Destheros doit aller au nord essayer de
# repérer si d'autres troupes andaries sont sur le sol nordien."
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("Il pense que Windalya va être attaquée.", $true, $false, $false, $false, $false, $true, 1, $false, "Il pense que Windalya va être attaquée, et décide de rentrer rapidos à Windalya, Kris avec lui. Katar Destheros doit aller au nord essayer de repérer si d’autres troupes andaries sont sur le sol nordien.", 2)
